# Auto-generated Excel COM-interop script
# Applies numeric corrections to currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets,
# matching the upstream scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3911.25
$ws.Range("I6").Value = 72.75
$ws.Range("J6").Value = 7749.75
$ws.Range("K6").Value = 218.25
$ws.Range("L6").Value = 23249.25
$ws.Range("M6").Value = -106.25
$ws.Range("N6").Value = -23473.25
$ws.Range("H19").Value = 15149.6875
$ws.Range("I19").Value = 2162.25
$ws.Range("J19").Value = 28137.125
$ws.Range("K19").Value = 2162.25
$ws.Range("L19").Value = 28137.125
$ws.Range("M19").Value = -1987.25
$ws.Range("N19").Value = -28487.125
$ws.Range("H70").Value = 1000759
$ws.Range("I70").Value = 3597077
$ws.Range("K70").Value = 10791231
$ws.Range("M70").Value = -10790961
$ws.Range("H73").Value = 1000759
$ws.Range("I73").Value = 3597077
$ws.Range("K73").Value = 10791231
$ws.Range("M73").Value = -10790295
$ws.Range("H103").Value = 735
$ws.Range("I103").Value = 628
$ws.Range("K103").Value = 1884
$ws.Range("M103").Value = -1298
$ws.Range("H111").Value = 32812.855
$ws.Range("I111").Value = 1156.75
$ws.Range("K111").Value = 3470.25
$ws.Range("M111").Value = -403.25
$ws.Range("H112").Value = 54171.105
$ws.Range("J112").Value = 1770.3636
$ws.Range("L112").Value = 5311.0908
$ws.Range("N112").Value = -7527.0908
$ws.Range("H132").Value = 4123.9062
$ws.Range("I132").Value = 3758.4443
$ws.Range("K132").Value = 11275.3329
$ws.Range("M132").Value = -8745.332900000001
$ws.Range("H137").Value = 1984.9524
$ws.Range("I137").Value = 1614.5385
$ws.Range("K137").Value = 4843.6155
$ws.Range("M137").Value = -2293.6155
$ws.Range("H138").Value = 2394.6758
$ws.Range("I138").Value = 1566.375
$ws.Range("K138").Value = 4699.125
$ws.Range("M138").Value = 440.875
$ws.Range("H141").Value = 841.8
$ws.Range("I141").Value = 934.6667
$ws.Range("K141").Value = 2804.0001
$ws.Range("M141").Value = 2375.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 37040150
$ws.Range("I61").Value = 41669044
$ws.Range("K61").Value = 41669044
$ws.Range("M61").Value = -41668832
$ws.Range("H74").Value = 1485.8077
$ws.Range("I74").Value = 1031
$ws.Range("K74").Value = 1031
$ws.Range("M74").Value = -157
$ws.Range("H77").Value = 1485.8077
$ws.Range("I77").Value = 1031
$ws.Range("K77").Value = 5155
$ws.Range("M77").Value = -787
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 10755285
$ws.Range("I122").Value = 11496822
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 34490466
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -34488016
$ws.Range("N122").Value = -13898.5
$ws.Range("I132").Value = 50001908
$ws.Range("K132").Value = 150005724
$ws.Range("M132").Value = -150003194
$ws.Range("H136").Value = 37040150
$ws.Range("I136").Value = 41669044
$ws.Range("K136").Value = 125007132
$ws.Range("M136").Value = -125004582

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 443.125
$ws.Range("I22").Value = 443.125
$ws.Range("K22").Value = 443.125
$ws.Range("M22").Value = -270.125
$ws.Range("H105").Value = 2100
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 2200
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 2200
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -5694
$ws.Range("H107").Value = 26562166
$ws.Range("I107").Value = 290619.06
$ws.Range("K107").Value = 290619.06
$ws.Range("M107").Value = -288699.06
$ws.Range("H134").Value = 3106.7083
$ws.Range("I134").Value = 3036
$ws.Range("K134").Value = 9108
$ws.Range("M134").Value = -6573

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1623.4
$ws.Range("I16").Value = 1582.6666
$ws.Range("K16").Value = 1582.6666
$ws.Range("M16").Value = -1295.6666
$ws.Range("H113").Value = 1623.4
$ws.Range("I113").Value = 1582.6666
$ws.Range("K113").Value = 1582.6666
$ws.Range("M113").Value = 587.3334
$ws.Range("H132").Value = 2891.5483
$ws.Range("I132").Value = 2897.2273
$ws.Range("K132").Value = 8691.6819
$ws.Range("M132").Value = -6161.6819
$ws.Range("H134").Value = 2888.8
$ws.Range("I134").Value = 2500
$ws.Range("J134").Value = 2986
$ws.Range("K134").Value = 7500
$ws.Range("L134").Value = 8958
$ws.Range("M134").Value = -4965
$ws.Range("N134").Value = -14028

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 695.4400000000001
$ws.Range("J113").Value = 879.4375
$ws.Range("L113").Value = 2638.3125
$ws.Range("N113").Value = -6978.3125
$ws.Range("H122").Value = 767.2857
$ws.Range("J122").Value = 876.3333
$ws.Range("L122").Value = 7886.9997
$ws.Range("N122").Value = -12786.9997
$ws.Range("H132").Value = 1453.6364
$ws.Range("I132").Value = 1453.6364
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13082.7276
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10552.7276
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1627.0962
$ws.Range("I102").Value = 1198.0555
$ws.Range("K102").Value = 1198.0555
$ws.Range("M102").Value = 423.9445000000001
$ws.Range("H107").Value = 1890.2
$ws.Range("J107").Value = 2268.5454
$ws.Range("L107").Value = 2268.5454
$ws.Range("N107").Value = -6108.5454
$ws.Range("H113").Value = 2426.0833
$ws.Range("J113").Value = 2901.4443
$ws.Range("L113").Value = 2901.4443
$ws.Range("N113").Value = -7241.4443
$ws.Range("H122").Value = 2246.2222
$ws.Range("I122").Value = 2009.5
$ws.Range("J122").Value = 3074.75
$ws.Range("K122").Value = 6028.5
$ws.Range("L122").Value = 9224.25
$ws.Range("M122").Value = -3578.5
$ws.Range("N122").Value = -14124.25
$ws.Range("H126").Value = 11818.363
$ws.Range("I126").Value = 19034.166
$ws.Range("K126").Value = 57102.49800000001
$ws.Range("M126").Value = -54632.49800000001
$ws.Range("H132").Value = 3768.55
$ws.Range("I132").Value = 3210.6875
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 9632.0625
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -7102.0625
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 71430560
$ws.Range("I7").Value = 83334820
$ws.Range("K7").Value = 83334820
$ws.Range("M7").Value = -83334708
$ws.Range("H40").Value = 3179.0908
$ws.Range("I40").Value = 2710.7144
$ws.Range("K40").Value = 2710.7144
$ws.Range("M40").Value = -2574.7144
$ws.Range("H122").Value = 3609.5833
$ws.Range("I122").Value = 2889.375
$ws.Range("K122").Value = 8668.125
$ws.Range("M122").Value = -6218.125
$ws.Range("H126").Value = 71430560
$ws.Range("I126").Value = 83334820
$ws.Range("K126").Value = 250004460
$ws.Range("M126").Value = -250001990
$ws.Range("H132").Value = 3520.4092
$ws.Range("I132").Value = 3452.45
$ws.Range("K132").Value = 10357.35
$ws.Range("M132").Value = -7827.349999999999
$ws.Range("H133").Value = 113517.8
$ws.Range("J133").Value = 113517.8
$ws.Range("L133").Value = 113517.8
$ws.Range("N133").Value = -118577.8
$ws.Range("H136").Value = 2201.4707
$ws.Range("I136").Value = 1995
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 5985
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -3435
$ws.Range("N136").Value = -16350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6666.3335
$ws.Range("I122").Value = 7500
$ws.Range("J122").Value = 4999
$ws.Range("K122").Value = 22500
$ws.Range("L122").Value = 14997
$ws.Range("M122").Value = -20050
$ws.Range("N122").Value = -19897
$ws.Range("H126").Value = 2062.375
$ws.Range("I126").Value = 1916.5
$ws.Range("K126").Value = 5749.5
$ws.Range("M126").Value = -3279.5
$ws.Range("H132").Value = 3226.9429
$ws.Range("I132").Value = 3131.6333
$ws.Range("K132").Value = 9394.8999
$ws.Range("M132").Value = -6864.8999
$ws.Range("H136").Value = 4786.647
$ws.Range("J136").Value = 6384.8887
$ws.Range("L136").Value = 19154.6661
$ws.Range("N136").Value = -24254.6661
